# Refresh FlashScore betting-odds values on Sheet1 (rows 2, 3, 4 and 7)
# to match the 2024-10-15 data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.7
$ws.Range("I2").Value = 5.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("X2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AH2").Value = 26
$ws.Range("AI2").Value = 19
$ws.Range("AJ2").Value = 51
$ws.Range("AQ2").Value = 34
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 29

# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.88

# Row 4
$ws.Range("G4").Value = 2.47
$ws.Range("H4").Value = 2.95
$ws.Range("I4").Value = 2.82
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.4
$ws.Range("S4").Value = 1.45
$ws.Range("T4").Value = 2.57
$ws.Range("U4").Value = 1.83
$ws.Range("W4").Value = 7.2
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 28
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 35
$ws.Range("AH4").Value = 14
$ws.Range("AI4").Value = 10.25
$ws.Range("AJ4").Value = 35
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 35
$ws.Range("AN4").Value = 4.4
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 22
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 300
$ws.Range("AT4").Value = 2.57
$ws.Range("AU4").Value = 6.9
$ws.Range("AV4").Value = 65
$ws.Range("AW4").Value = 4.75
$ws.Range("AX4").Value = 15.5
$ws.Range("AZ4").Value = 75

# Row 7
$ws.Range("G7").Value = 5.7
$ws.Range("H7").Value = 3.7
$ws.Range("I7").Value = 1.53
$ws.Range("J7").Value = 5.8
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.12
$ws.Range("N7").Value = 7.2
$ws.Range("O7").Value = 1.31
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.82
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.67
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.72
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 7.4
$ws.Range("AF7").Value = 100
$ws.Range("AG7").Value = 6.1
$ws.Range("AL7").Value = 30
$ws.Range("AM7").Value = 900
$ws.Range("AN7").Value = 7.2
$ws.Range("AP7").Value = 40
$ws.Range("AT7").Value = 2.67
$ws.Range("AV7").Value = 90
$ws.Range("AX7").Value = 7.5
$ws.Range("AY7").Value = 19
$ws.Range("AZ7").Value = 25
$ws.Range("BA7").Value = 65
$ws.Range("BB7").Value = 300
